{"js": "// Add a new \"Navigation for blind people\" use case, inserted right\n// before the existing \"Environmental Study:\" list item.\n//\n// The new use case consists of two paragraphs:\n//   1. A bulleted/numbered heading paragraph: \"Navigation for blind people\"\n//      (same list + style as the other use-case headings, e.g. \"Environmental\n//      Study:\", \"Robotics and Navigation:\", ...).\n//   2. A normal (non-list) body paragraph with the description text.\n\nconst body = context.document.body;\n\n// Locate the \"Environmental Study:\" heading paragraph - the new use case is\n// inserted directly before it.\nconst searchResults = body.search(\"Environmental Study:\", { matchCase: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nconst anchorHeading = searchResults.items[0].paragraphs.getFirst();\nanchorHeading.load(\"style\");\nanchorHeading.list.load(\"id\");\nawait context.sync();\n\n// The paragraph right after the heading is the existing \"Normal\"-styled\n// body paragraph (\"Ecologists can analyze ...\") - reuse its style for the\n// new body paragraph so the new use case matches the formatting of the\n// others.\nconst anchorBody = anchorHeading.getNext();\nanchorBody.load(\"style\");\nawait context.sync();\n\nconst headingStyle = anchorHeading.style;\nconst bodyStyle = anchorBody.style;\nconst listId = anchorHeading.list.id;\n\n// Insert the new heading paragraph before \"Environmental Study:\".\nconst newHeading = anchorHeading.insertParagraph(\n  \"Navigation for blind people\",\n  Word.InsertLocation.before\n);\nnewHeading.style = headingStyle;\nawait context.sync();\n\n// Attach the new heading to the same numbered list (level 0) as the other\n// use-case headings.\nnewHeading.attachToList(listId, 0);\nawait context.sync();\n\n// Insert the new body paragraph right after the new heading.\nconst bodyText =\n  \"Blind individuals can enhance their navigation by gaining awareness of their surroundings by using Lidar scanning. Instead of using traditional methods like canes that heavily rely on tactile feedback and can sometimes be imprecise, LiDAR scanner assists them in identifying and measuring the distance to objects including locating features like doors. Real-Time Mapping and combination with other sensor such as audio sensor also help them  to navigate their surroundings more efficiently and safely. A danger alert can also be implemented so that they could prevent themselve from involve in any accident.\";\n\nconst newBody = newHeading.insertParagraph(bodyText, Word.InsertLocation.after);\nnewBody.style = bodyStyle;\nawait context.sync();\n", "ps1": "# Add a new \"Navigation for blind people\" use case, inserted right before\n# the existing \"Environmental Study:\" list item.\n#\n# The new use case consists of two paragraphs:\n#   1. A bulleted/numbered heading paragraph: \"Navigation for blind people\"\n#      (same list + style as the other use-case headings, e.g. \"Environmental\n#      Study:\", \"Robotics and Navigation:\", ...).\n#   2. A normal (non-list) body paragraph with the description text.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Environmental Study:\" heading paragraph - the new use case is\n# inserted directly before it.\n$searchRange = $d.Content\n$searchRange.Find.ClearFormatting()\n$searchRange.Find.Text = \"Environmental Study:\"\n$searchRange.Find.MatchCase = $false\n$searchRange.Find.MatchWholeWord = $false\n$searchRange.Find.MatchWildcards = $false\n$searchRange.Find.Execute() | Out-Null\n\n$headingAnchor = $searchRange.Paragraphs(1)\n\n# The paragraph right after the heading is the existing \"Normal\"-styled body\n# paragraph (\"Ecologists can analyze ...\") - reuse its style for the new body\n# paragraph so the new use case matches the formatting of the others.\n$bodyStyleRef = $headingAnchor.Next()\n$bodyStyleName = $bodyStyleRef.Style.NameLocal\n\n$bodyText = \"Blind individuals can enhance their navigation by gaining awareness of their surroundings by using Lidar scanning. Instead of using traditional methods like canes that heavily rely on tactile feedback and can sometimes be imprecise, LiDAR scanner assists them in identifying and measuring the distance to objects including locating features like doors. Real-Time Mapping and combination with other sensor such as audio sensor also help them  to navigate their surroundings more efficiently and safely. A danger alert can also be implemented so that they could prevent themselve from involve in any accident.\"\n\n# Insert the new heading paragraph and new body paragraph, both before\n# \"Environmental Study:\". Inserting plain text + paragraph marks this way\n# makes both new paragraphs inherit the heading's list style/numbering.\n$insertPoint = $headingAnchor.Range.Duplicate()\n$insertPoint.Collapse(1)\n$insertPoint.InsertBefore(\"Navigation for blind people`r\" + $bodyText + \"`r\")\n\n# Fix up the new body paragraph's style: it should be the plain \"Normal\"\n# style (no list numbering), not the inherited list style.\n$bodySearchRange = $d.Content\n$bodySearchRange.Find.ClearFormatting()\n$bodySearchRange.Find.Text = \"Blind individuals can enhance their navigation\"\n$bodySearchRange.Find.MatchCase = $false\n$bodySearchRange.Find.MatchWholeWord = $false\n$bodySearchRange.Find.MatchWildcards = $false\n$bodySearchRange.Find.Execute() | Out-Null\n\n$newBodyPara = $bodySearchRange.Paragraphs(1)\n$newBodyPara.Style = $bodyStyleName\n\n# NOTE: the new heading paragraph already carries the correct style and list\n# numbering, inherited from $headingAnchor via InsertBefore above - do not\n# re-assign its .Style here, since re-assigning a paragraph's Style to its\n# own current style name strips the direct list-numbering (numPr) format.\n"}
